# Updates the '想去人数' (want-to-go count) figures in column F across all four
# worksheets to reflect refreshed scrape output, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 7350
$ws.Range("F10").Value = 1928
$ws.Range("F11").Value = 5365
$ws.Range("F14").Value = 7415
$ws.Range("F15").Value = 8711
$ws.Range("F18").Value = 859
$ws.Range("F19").Value = 4361
$ws.Range("F26").Value = 86
$ws.Range("F27").Value = 1618
$ws.Range("F28").Value = 689
$ws.Range("F29").Value = 876
$ws.Range("F30").Value = 1847
$ws.Range("F35").Value = 1416
$ws.Range("F39").Value = 379
$ws.Range("F40").Value = 2909
$ws.Range("F41").Value = 4007
$ws.Range("F44").Value = 405
$ws.Range("F46").Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F26").Value = 98

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5089

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5089
$ws.Range("F10").Value = 5365
$ws.Range("F12").Value = 7415
$ws.Range("F16").Value = 859
$ws.Range("F17").Value = 4361
$ws.Range("F24").Value = 86
$ws.Range("F25").Value = 1618
$ws.Range("F26").Value = 689
$ws.Range("F27").Value = 876
$ws.Range("F28").Value = 1847
$ws.Range("F38").Value = 379
$ws.Range("F39").Value = 98
$ws.Range("F40").Value = 4007
$ws.Range("F44").Value = 405
